$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Copy formatting from suitable "donor" rows so the new rows inherit the
#    same cell styles as the existing data, without creating extra/duplicate
#    style records.
#    - Row 80   : WLAN/Laptop pattern (A s2,B s1,C s2,D s2,E s2,F s4,G s2)
#    - Row 106  : same as row 80 but E uses the "big number" style (s3)
#    - Row 236  : iPad/WLAN pattern (A s2,B s1,C/D/E no style,F s4,G no style)
# ---------------------------------------------------------------------------

$ws.Range("A80:G80").Copy()
$ws.Range("A256:G259").PasteSpecial(-4122) | Out-Null
$ws.Range("A261:G269").PasteSpecial(-4122) | Out-Null

$ws.Range("A106:G106").Copy()
$ws.Range("A260:G260").PasteSpecial(-4122) | Out-Null

$ws.Range("A80:G80").Copy()
$ws.Range("A270:G270").PasteSpecial(-4122) | Out-Null

$ws.Range("A236:G236").Copy()
$ws.Range("A271:G275").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Fill in the new row values (rows 256-275)
# ---------------------------------------------------------------------------

$rowsData = @(
    @(256, "24.01.2016", 0.66875000000000007, 11.12,               3.25,  38,    "WLAN", "Laptop"),
    @(257, "24.01.2016", 0.67152777777777783, 11.9,                1.1000000000000001, 27, "WLAN", "Laptop"),
    @(258, "24.01.2016", 0.67291666666666661, 10.11,               12.05, 32,    "WLAN", "Laptop"),
    @(259, "24.01.2016", 0.6743055555555556,  11.23,               6.89,  117,   "WLAN", "Laptop"),
    @(260, "24.01.2016", 0.67569444444444438, 11.73,               25.73, 11638, "WLAN", "Laptop"),
    @(261, "27.01.2016", 0.86736111111111114, 9.07,                15.66, 578,   "WLAN", "Laptop"),
    @(262, "27.01.2016", 0.8666666666666667,  3.73,                1.19,  42,    "WLAN", "Laptop"),
    @(263, "27.01.2016", 0.86597222222222225, 9.09,                14.62, 27,    "WLAN", "Laptop"),
    @(264, "27.01.2016", 0.86458333333333337, 9.19,                1.06,  23,    "WLAN", "Laptop"),
    @(265, "27.01.2016", 0.5,                 9.1,                 6.37,  27,    "WLAN", "Laptop"),
    @(266, "27.01.2016", 0.93263888888888891, 11.64,               23.05, 27,    "WLAN", "Laptop"),
    @(267, "27.01.2016", 0.93055555555555547, 11.17,               1.08,  28,    "WLAN", "Laptop"),
    @(268, "27.01.2016", 0.92569444444444438, 7.48,                17,    31,    "WLAN", "Laptop"),
    @(269, "27.01.2016", 0.92361111111111116, 11.24,               20.73, 44,    "WLAN", "Laptop"),
    @(270, "27.01.2016", 0.92083333333333339, 10.82,               1.77,  25,    "WLAN", "Laptop"),
    @(271, "28.01.2016", 0.37222222222222223, 9.42,                0.99,  37,    "WLAN", "iPad"),
    @(272, "28.01.2016", 0.37152777777777773, 8.9600000000000009,  0.82,  38,    "WLAN", "iPad"),
    @(273, "28.01.2016", 0.36944444444444446, 3.26,                0.81,  38,    "WLAN", "iPad"),
    @(274, "28.01.2016", 0.36944444444444446, 6.19,                1.04,  40,    "WLAN", "iPad"),
    @(275, "28.01.2016", 0.36874999999999997, 7.81,                0.7,   45,    "WLAN", "iPad")
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 5).Value2 = $row[5]
    $ws.Cells.Item($r, 6).Value2 = $row[6]
    $ws.Cells.Item($r, 7).Value2 = $row[7]
    $ws.Cells.Item($r, 8).Value2 = 15.747999999999999
    $ws.Cells.Item($r, 9).Value2 = 1.1020000000000001
}

# ---------------------------------------------------------------------------
# 3) Row 270's H/I cells get a new explicit style (font color forced to
#    automatic/black -> introduces a new cellXfs entry, matching the diff).
# ---------------------------------------------------------------------------

$ws.Range("H270:I270").Font.Color = 0

# ---------------------------------------------------------------------------
# 4) Update the sheet view (scroll position / selection) to match the new
#    bottom of the data.
# ---------------------------------------------------------------------------

$ws.Activate()
$ws.Range("C271:E275").Select()
